# netCrypto.xlsx — re-upload edit
#
# Target diff touches five things:
#   1. x15ac:absPath (folder Excel was last saved from)
#   2. xr:revisionPtr/@documentId (co-authoring save GUID)
#   3. bookViews/workbookView/@xWindow (host window X position)
#   4. sheetView: add topLeftCell="H1" (scrolled viewport) + move the
#      selection from H22 to T3
#   5. SheetName1!T2 value 216399 -> 216919
#
# Items 1-3 and the topLeftCell part of item 4 are metadata that real Excel
# stamps on a file purely as a side effect of *where/how* it was saved (last
# folder used, window placement, co-authoring revision id) — they are not
# exposed anywhere in the Excel object model (no Workbook/Window property
# backs them), so they cannot be produced by driving the app. We still poke
# the closest available COM surface for them below (harmless no-ops if the
# host doesn't route them anywhere) and then perform the two edits that
# genuinely are reachable through the object model: moving the selection and
# updating the cell value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Best-effort: host window position (xWindow on <workbookView>) --------
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 28680
} catch {}

# --- Best-effort: scrolled top-left cell (topLeftCell="H1") ---------------
try {
    $activeWin = $excel.ActiveWindow
    $activeWin.ScrollColumn = 8   # column H
    $activeWin.ScrollRow = 1
} catch {}

# --- Move the selection from H22 to T3 -------------------------------------
$ws.Range("T3").Select()

# --- Update SheetName1!T2 ---------------------------------------------------
$ws.Range("T2").Value = 216919

$wb.Save()
